# Weekly update: insert a new price record as the new row 20 (date 44603),
# pushing all the subsequent daily records down by one row. The last
# existing row (29) becomes the new row 30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. First, create the brand-new last row (30) as a copy of the
#        current last row (29), before row 29's data gets overwritten. ---
$ws.Cells.Item(30, 1).Value = $ws.Cells.Item(29, 1).Value2
$ws.Cells.Item(30, 2).Value = $ws.Cells.Item(29, 2).Value2
$ws.Cells.Item(30, 3).Value = $ws.Cells.Item(29, 3).Value2
$ws.Cells.Item(30, 4).Value = $ws.Cells.Item(29, 4).Value2
$ws.Cells.Item(30, 4).NumberFormat = $ws.Cells.Item(29, 4).NumberFormat
$ws.Cells.Item(30, 5).Value = $ws.Cells.Item(29, 5).Value2
$ws.Cells.Item(30, 6).Value = $ws.Cells.Item(29, 6).Value2
$ws.Cells.Item(30, 7).Value = $ws.Cells.Item(29, 7).Value2
$ws.Cells.Item(30, 8).Value = $ws.Cells.Item(29, 8).Value2
$ws.Cells.Item(30, 9).Value = $ws.Cells.Item(29, 9).Value2
$ws.Cells.Item(30, 10).Value = $ws.Cells.Item(29, 10).Value2
$ws.Cells.Item(30, 11).Value = $ws.Cells.Item(29, 11).Value2
$ws.Cells.Item(30, 12).Value = $ws.Cells.Item(29, 12).Value2
$ws.Cells.Item(30, 13).Value = $ws.Cells.Item(29, 13).Value2
$ws.Cells.Item(30, 14).Value = $ws.Cells.Item(29, 14).Value2
$ws.Cells.Item(30, 15).Value = $ws.Cells.Item(29, 15).Value2
$ws.Cells.Item(30, 16).Value = $ws.Cells.Item(29, 16).Value2
$ws.Cells.Item(30, 17).Value = $ws.Cells.Item(29, 17).Value2
$ws.Cells.Item(30, 18).Value = $ws.Cells.Item(29, 18).Value2

# --- 2. Shift the varying columns (D, J, K, L, M, P) of rows 20..28 down
#        into rows 21..29. Walk bottom-up so a source row is never
#        overwritten before it has been read. All other columns
#        (A,B,C,E,F,G,H,I,N,O,Q,R) are identical for every data row,
#        so they need no change. ---
for ($r = 29; $r -ge 21; $r--) {
    $src = $r - 1
    $ws.Cells.Item($r, 4).Value = $ws.Cells.Item($src, 4).Value2
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($src, 10).Value2
    $ws.Cells.Item($r, 11).Value = $ws.Cells.Item($src, 11).Value2
    $ws.Cells.Item($r, 12).Value = $ws.Cells.Item($src, 12).Value2
    $ws.Cells.Item($r, 13).Value = $ws.Cells.Item($src, 13).Value2
    $ws.Cells.Item($r, 16).Value = $ws.Cells.Item($src, 16).Value2
}

# --- 3. Finally, write the brand-new record into row 20. ---
$ws.Cells.Item(20, 4).Value = 44603
$ws.Cells.Item(20, 10).Value = 250
$ws.Cells.Item(20, 11).Value = 2500
$ws.Cells.Item(20, 12).Value = 3000
$ws.Cells.Item(20, 13).Value = 2750
$ws.Cells.Item(20, 16).Value = 2750
